$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")
$ws.Activate()

# "cambio de fracciones e historico": roll the reporting period / history
# dates forward from Q2 2022 to Q3 2022.
$ws.Range("B8").Value = 44743   # Fecha de inicio del periodo que se informa (01/07/2022)
$ws.Range("C8").Value = 44834   # Fecha de termino del periodo que se informa (30/09/2022)
$ws.Range("K8").Value = 44844   # Fecha de validacion (10/10/2022)
$ws.Range("L8").Value = 44844   # Fecha de actualizacion (10/10/2022)

# Scroll the sheet back to the left (was parked at column K) and move the
# active cell/selection to C11, matching the saved workbook's view state.
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("C11").Select()
